$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sec invoice Master")

# --- New secondary-invoice batch (qq / custom suite update) ---------------
# Three invoices were processed in this batch; the master row (row 2) is
# refreshed to reflect the last one of the three:
#   58327782 / FCT915497466927775744 / 58327782+1 / 105.22
#   58327789 / FCT915506118178897920 / 58327789+1 / 85.1
#   58327790 / FCT915518758829686784 / 58327790+1 / 112.84   <- lands in row 2

$orderIds  = @("58327782", "58327789", "58327790")
$tracking  = @("FCT915497466927775744", "FCT915506118178897920", "FCT915518758829686784")
$invoiceNo = @("58327782+1", "58327789+1", "58327790+1")
$amounts   = @("105.22", "85.1", "112.84")

# A "scratch" cell used purely to funnel numeric-looking text (order id /
# amount) through Excel as genuine text rather than a number, so the cell
# keeps its original General number format / style family (same as the
# pre-existing row) instead of picking up a new "@" text format.
$scratch = $ws.Range("ZZ1")

for ($i = 0; $i -lt $orderIds.Length; $i++) {

    # FC Order ID (A2) -- numeric-looking, must stay text
    $scratch.NumberFormat = "@"
    $scratch.Value = $orderIds[$i]
    $scratch.Copy()
    $ws.Range("A2").PasteSpecial(-4163) | Out-Null
    $scratch.Clear()

    # Tracking # (C2) -- already non-numeric, plain assignment keeps text
    $ws.Range("C2").Value = $tracking[$i]

    # New Invoice Amount (F2) -- numeric-looking, must stay text
    $scratch.NumberFormat = "@"
    $scratch.Value = $amounts[$i]
    $scratch.Copy()
    $ws.Range("F2").PasteSpecial(-4163) | Out-Null
    $scratch.Clear()

    # SECONDARY INV # (I2) -- already non-numeric, plain assignment keeps text
    $ws.Range("I2").Value = $invoiceNo[$i]
}

$excel.CutCopyMode = 0
